$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing "sum" header (G1) onto the new "Save" header (H1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new header text and value
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
